$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E9, E10, E12, E16 switch from "NO" to "SI" (login requerido column)
$ws.Range("E9").Value = "SI"
$ws.Range("E10").Value = "SI"
$ws.Range("E12").Value = "SI"
$ws.Range("E16").Value = "SI"

# Row 34 - Caso #28: Crear Blog
$ws.Range("A34").Value = "Caso #28"
$ws.Range("B34").Value = "Crear Blog"
$ws.Range("C34").Value = "Crear Blog desde el menu propio de admin"
$ws.Range("D34").Value = 45276
$ws.Range("E34").Value = "SI"
$ws.Range("F34").Value = "-"
$ws.Range("G34").Value = "OK"

# Row 35 - Caso #29: Ver listado de blogs (admin)
$ws.Range("A35").Value = "Caso #29"
$ws.Range("B35").Value = "Ver listado de blogs (admin)"
$ws.Range("C35").Value = "Ver listado de Blogs desde el menu propio de admin"
$ws.Range("D35").Value = 45276
$ws.Range("E35").Value = "SI"
$ws.Range("F35").Value = "-"
$ws.Range("G35").Value = "OK"

# Update the active selection to match the last edit made
$ws.Range("E35:G35").Select()
